$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper pattern -------------------------------------------------
# For cells that must keep TEXT type even though the literal looks like a
# number (e.g. "001010", "3.6000", "4", "0.25" ...), we assign the value
# with a leading apostrophe (forces Excel to treat it as text) and then
# reset the cell Style back to "Normal" so no stray number-format style
# is left behind on the cell (keeps XML style index identical to before).

# ---------------------------------------------------------------------
# Row 1 headers
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Numbers"
$ws.Range("C1").Value = "Binary_numbers"
$ws.Range("D1").Value = "Values_x"
$ws.Range("E1").Value = "Values_f(x)"

# ---------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = "'001010"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 3.6
$ws.Range("E2").Value = "'3.6000"
$ws.Range("E2").Style = "Normal"

# ---------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = "'101000"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 5.4
$ws.Range("E3").Value = "'5.4000"
$ws.Range("E3").Style = "Normal"

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = "'001101"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 3.78
$ws.Range("E4").Value = "'3.7800"
$ws.Range("E4").Style = "Normal"

# ---------------------------------------------------------------------
# Row 5 (B5, C5, D5 unchanged; only E5 changes)
# ---------------------------------------------------------------------
$ws.Range("E5").Value = "'5.5200"
$ws.Range("E5").Style = "Normal"

# ---------------------------------------------------------------------
# Row 7 header: drop "prob_of_crossing" (old E7), shift F7:O7 -> E7:N7,
# clear O7.
# ---------------------------------------------------------------------
$ws.Range("E7").Value = "prob_mut_ind"
$ws.Range("F7").Value = "prob_mut_gen"
$ws.Range("G7").Value = "function_entry"
$ws.Range("H7").Value = "delta_x"
$ws.Range("I7").Value = "find_x_by"
$ws.Range("J7").Value = "iterator_entry"
$ws.Range("K7").Value = "range"
$ws.Range("L7").Value = "jump_numbers"
$ws.Range("M7").Value = "points_numbers"
$ws.Range("N7").Value = "bits_required"
$ws.Range("O7").Clear()

# ---------------------------------------------------------------------
# Row 8: drop old E8 (0.75), shift F8:O8 -> E8:N8. The old H8 "x**2"
# (string) becomes new G8 "x" (string); old I8 0.06 (number) becomes new
# H8 0.06 (number); etc.
# ---------------------------------------------------------------------
$ws.Range("E8").Value = "'0.25"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'0.35"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'x"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 0.06
$ws.Range("I8").Value = "Minimización"
$ws.Range("J8").Value = "'2"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 33.3333
$ws.Range("M8").Value = 34.3333
$ws.Range("N8").Value = 6
$ws.Range("O8").Clear()
